# Apply data-refresh edits to the "TPE" results table on Feuil1,
# then restore the selection/scroll position captured in the source file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")
$ws.Activate()

# Row 36 (T26_3 family, Coarse)
$ws.Range("D36").Value = 455
$ws.Range("E36").Value = 2849
$ws.Range("F36").Value = 24594

# Row 43 (T15_8 family, Coarse)
$ws.Range("D43").Value = 988
$ws.Range("E43").Value = 1063
$ws.Range("F43").Value = 23892

# Row 44 (T15_8 family, 25)
$ws.Range("D44").Value = 3621
$ws.Range("E44").Value = 3762
$ws.Range("F44").Value = 24046

# Row 45 (T15_8 family, 20)
$ws.Range("D45").Value = 5653
$ws.Range("E45").Value = 5829
$ws.Range("F45").Value = 24070

# Row 50 (T26_6 family, Coarse)
$ws.Range("D50").Value = 988
$ws.Range("E50").Value = 4103
$ws.Range("F50").Value = 24572

# Row 51 (T26_6 family, 25)
$ws.Range("D51").Value = 3621
$ws.Range("E51").Value = 14767
$ws.Range("F51").Value = 24296

# Row 52 (T26_6 family, 20)
$ws.Range("D52").Value = 5653
$ws.Range("E52").Value = 22965
$ws.Range("F52").Value = 24234

# Row 54 (T26_6 family, 5)
$ws.Range("F54").Value = 24238

# Restore the saved view state: selection on F54
$ws.Range("F54").Select()
